$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 376 (pushes the existing rows 376:401 down to 377:402,
# and extends the used range to A1:R402).
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(376, 1).Value = 8
$ws.Cells.Item(376, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(376, 3).Value = "Coquimbo"
$ws.Cells.Item(376, 4).Value = 44826
$ws.Cells.Item(376, 5).Value = 4
$ws.Cells.Item(376, 6).Value = 100112032
$ws.Cells.Item(376, 7).Value = "Zapallo italiano"
$ws.Cells.Item(376, 8).Value = "Sin especificar"
$ws.Cells.Item(376, 9).Value = "Primera"
$ws.Cells.Item(376, 10).Value = 600
$ws.Cells.Item(376, 11).Value = 11000
$ws.Cells.Item(376, 12).Value = 12000
$ws.Cells.Item(376, 13).Value = 11500
$ws.Cells.Item(376, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(376, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(376, 16).Value = 230
$ws.Cells.Item(376, 17).Value = 50
$ws.Cells.Item(376, 18).Value = "Hortaliza"
